$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.232.79"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "1.788.30"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "225.94"
$ws.Range("E5").Value = "  -0.77%  "
$ws.Range("D6").Value = "0.556"
$ws.Range("E6").Value = "  +1.59%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "32.27"
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").Value = "0.0688"
$ws.Range("E10").Value = "  -0.39%  "
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("D12").Value = "2.047.42"
$ws.Range("D13").Value = "1.803.15"
$ws.Range("E13").Value = "  +0.62%  "
$ws.Range("D14").Value = "11.01"
$ws.Range("E14").Value = "  -4.87%  "
$ws.Range("D15").Value = "34.227.59"
$ws.Range("D16").Value = "0.624"
$ws.Range("E16").Value = "  +0.16%  "
$ws.Range("D17").Value = "4.19"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").Value = "67.97"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "245.87"
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0798"
$ws.Range("E20").Value = "  +2.05%  "
$ws.Range("D21").Value = "10.93"
$ws.Range("E21").Value = "  +0.51%  "
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").Value = "4.14"
$ws.Range("E23").Value = "  +0.94%  "
$ws.Range("E24").Value = "  +0.36%  "
$ws.Range("D25").Value = "161.51"
$ws.Range("E25").Value = "  -0.42%  "
$ws.Range("D26").Value = "7.16"
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").Value = "16.33"
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("E28").Value = "  +0.97%  "
$ws.Range("D29").Value = "1.01"
$ws.Range("E29").Value = "  +0.37%  "
$ws.Range("E30").Value = "  -1.01%  "
$ws.Range("D31").Value = "0.0520"
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("E32").Value = "  +2.20%  "
$ws.Range("D33").Value = "3.75"
$ws.Range("E33").Value = "  +3.66%  "
$ws.Range("D34").Value = "1.80"
$ws.Range("E34").Value = "  -1.81%  "
$ws.Range("D35").Value = "1.437.89"
$ws.Range("E35").Value = "  +1.19%  "
$ws.Range("D36").Value = "2.58"
$ws.Range("E36").Value = "  +9.23%  "
$ws.Range("E37").Value = "  +3.40%  "
$ws.Range("E38").Value = "  +1.30%  "
$ws.Range("E39").Value = "  -0.99%  "
$ws.Range("D40").Value = "81.85"
$ws.Range("E40").Value = "  +1.45%  "
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "14.08"
$ws.Range("E41").Value = "  +5.62%  "
$ws.Range("B42").Value = "HuobiToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D42").Value = "2.39"
$ws.Range("E42").Value = "  +1.29%  "
$ws.Range("D43").Value = "2.73"
$ws.Range("E43").Value = "  +1.40%  "
$ws.Range("D44").Value = "0.919"
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("D45").Value = "0.0520"
$ws.Range("E45").Value = "  +2.28%  "
$ws.Range("D46").Value = "6.07"
$ws.Range("E46").Value = "  +0.14%  "
$ws.Range("E47").Value = "  +0.31%  "
$ws.Range("D48").Value = "1.943.96"
$ws.Range("E48").Value = "  -0.27%  "
$ws.Range("D49").Value = "105.33"
$ws.Range("E49").Value = "  -2.09%  "
$ws.Range("E50").Value = "  +0.13%  "
$ws.Range("E51").Value = "  -6.25%  "
